{"js": "// Adds a new list paragraph (ilvl=1, numId=1) at the very end of the\n// document body, right after the last existing paragraph\n// (\"...Tem que separar pela barra, caso contr\u00e1rio n\u00e3o funciona.\"),\n// explaining the :nth-child() selector usage.\n\nconst runsData = [\n  { t: \" Podemos selecionar um item espec\u00edfico de v\u00e1rios que possuem a mesma classe usando a extens\u00e3o :\", b: false, i: false },\n  { t: \"nth-child(n\u00famero){}\", b: false, i: true },\n  { t: \" logo ap\u00f3s o nome da classe que desejamos usar. Por\u00e9m com uma observa\u00e7\u00e3o: O n\u00famero colocado \u00e9 referente ao elemento/divis\u00e3o m\u00e3e desse filho, e n\u00e3o a\", b: false, i: false },\n  { t: \" quantidade de elementos com as mesmas classes, ou seja,\", b: false, i: false },\n  { t: \" se voc\u00ea possui v\u00e1rios elementos com a classe \", b: false, i: false },\n  { t: \".destaques__secundario\", b: false, i: true },\n  { t: \" que vem logo ap\u00f3s uma classe \", b: false, i: false },\n  { t: \".destaques__primario\", b: false, i: true },\n  { t: \", todas est\u00e3o dentro da divis\u00e3o com classe \", b: false, i: false },\n  { t: \".destaques\", b: false, i: true },\n  { t: \" e voc\u00ea quer fazer a configura\u00e7\u00e3o espec\u00edfica do primeiro elemento da classe \", b: false, i: false },\n  { t: \".destaques__secundario\", b: false, i: true },\n  { t: \" o seu \", b: false, i: false },\n  { t: \":nth-child(n\u00famero){}\", b: false, i: true },\n  { t: \" ter\u00e1 o n\u00famero 2, pois antes dele vem apenas 1 elemento da classe \", b: false, i: false },\n  { t: \".destaques__primario\", b: false, i: true },\n  { t: \". Se tivessem 2 elementos prim\u00e1rios antes do primeiro secund\u00e1rio, o n\u00famero do child do primeiro secund\u00e1rio seria \", b: false, i: false },\n  { t: \"3\", b: true, i: true },\n  { t: \" ao inv\u00e9s de 2\", b: false, i: false },\n  { t: \".\", b: false, i: false },\n];\n\n// Escape text for safe inclusion inside OOXML <w:t> elements.\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\")\n    .replace(/'/g, \"&apos;\");\n}\n\nfunction buildRunXml(r) {\n  const parts = [\n    '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>',\n  ];\n  if (r.b) {\n    parts.push(\"<w:b/><w:bCs/>\");\n  }\n  if (r.i) {\n    parts.push(\"<w:i/><w:iCs/>\");\n  }\n  parts.push('<w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/>');\n  const rPr = \"<w:rPr>\" + parts.join(\"\") + \"</w:rPr>\";\n  const text = escapeXml(r.t);\n  const preserve = text !== text.trim() || text.length === 0;\n  const t = preserve\n    ? '<w:t xml:space=\"preserve\">' + text + \"</w:t>\"\n    : \"<w:t>\" + text + \"</w:t>\";\n  return \"<w:r>\" + rPr + t + \"</w:r>\";\n}\n\nconst runsXml = runsData.map(buildRunXml).join(\"\");\n\nconst paragraphXml =\n  \"<w:p>\" +\n  \"<w:pPr>\" +\n  '<w:pStyle w:val=\"PargrafodaLista\"/>' +\n  \"<w:numPr>\" +\n  '<w:ilvl w:val=\"1\"/>' +\n  '<w:numId w:val=\"1\"/>' +\n  \"</w:numPr>\" +\n  '<w:jc w:val=\"both\"/>' +\n  \"<w:rPr>\" +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  \"<w:b/><w:bCs/>\" +\n  '<w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/>' +\n  '<w:u w:val=\"single\"/>' +\n  \"</w:rPr>\" +\n  \"</w:pPr>\" +\n  runsXml +\n  \"</w:p>\";\n\nconst ooxmlPackage =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  paragraphXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst insertionRange = lastParagraph.getRange(Word.RangeLocation.after);\ninsertionRange.insertOoxml(ooxmlPackage, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Adds a new list paragraph (ilvl=1, numId=1) at the very end of the\n# document body, right after the last existing paragraph\n# (\"...Tem que separar pela barra, caso contr\u00e1rio n\u00e3o funciona.\"),\n# explaining the :nth-child() selector usage.\n\n$d = $word.ActiveDocument\n\n# Build the OOXML for the new paragraph (WordprocessingML wrapped in the\n# standard \"single XML part\" package format expected by Range.InsertXML).\n$paragraphXml = '<w:p><w:pPr><w:pStyle w:val=\"PargrafodaLista\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr><w:jc w:val=\"both\"/><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:b/><w:bCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:u w:val=\"single\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> Podemos selecionar um item espec\u00edfico de v\u00e1rios que possuem a mesma classe usando a extens\u00e3o :</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:i/><w:iCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>nth-child(n\u00famero){}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> logo ap\u00f3s o nome da classe que desejamos usar. Por\u00e9m com uma observa\u00e7\u00e3o: O n\u00famero colocado \u00e9 referente ao elemento/divis\u00e3o m\u00e3e desse filho, e n\u00e3o a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> quantidade de elementos com as mesmas classes, ou seja,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> se voc\u00ea possui v\u00e1rios elementos com a classe </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:i/><w:iCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>.destaques__secundario</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> que vem logo ap\u00f3s uma classe </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:i/><w:iCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>.destaques__primario</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">, todas est\u00e3o dentro da divis\u00e3o com classe </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:i/><w:iCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>.destaques</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> e voc\u00ea quer fazer a configura\u00e7\u00e3o espec\u00edfica do primeiro elemento da classe </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:i/><w:iCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>.destaques__secundario</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> o seu </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:i/><w:iCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>:nth-child(n\u00famero){}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> ter\u00e1 o n\u00famero 2, pois antes dele vem apenas 1 elemento da classe </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:i/><w:iCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>.destaques__primario</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">. Se tivessem 2 elementos prim\u00e1rios antes do primeiro secund\u00e1rio, o n\u00famero do child do primeiro secund\u00e1rio seria </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> ao inv\u00e9s de 2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r></w:p>'\n\n$packageXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">' + `\n    '<pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n    '<w:body>' + $paragraphXml + '</w:body>' + `\n    '</w:document>' + `\n    '</pkg:xmlData>' + `\n    '</pkg:part>' + `\n    '</pkg:package>'\n\n# Position an empty range at the very end of the document body (right\n# before the final section break) and insert the new paragraph there.\n$endPos = $d.Content.End\n$insertionRange = $d.Range($endPos, $endPos)\n$insertionRange.InsertXML($packageXml) | Out-Null\n"}
